# Apply "repull data, push all data, mean calculation" edit:
# Column F (dSF) values were re-pulled and differ from column E (dS0) values
# for several rows. Update the affected cells on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    4  = -3
    11 = -4
    13 = -2
    15 = 3
    19 = -3
    20 = -9
    21 = -7
    22 = 3
    23 = -4
    24 = -3
    25 = -3
    26 = -1
    27 = -1
    29 = -4
    30 = -6
    31 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
